$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: A7 changes from "{{312312}}{{123123}}" to "dsadsa";
# C7 and D7 get new text values ("1231313" and "dsadqweqeasd").
$ws.Range("A7").Value = "dsadsa"

# "1231313" looks like a pure number -- prefix with a single quote so Excel
# stores it as text (shared string) instead of coercing it to a numeric value.
$ws.Range("C7").Value = "'1231313"

$ws.Range("D7").Value = "dsadqweqeasd"

# Row 10 / Row 11: clear the old placeholder text, leaving a (quote-prefixed)
# empty text value rather than a fully blank cell.
$ws.Range("A10").Value = "'"
$ws.Range("D11").Value = "'"

# G11 gets a brand new text value.
$ws.Range("G11").Value = "fyw是傻逼"

# Update the active selection to F4:F6 (a merged cell).
$ws.Range("F4:F6").Select()
